$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3686.4443
$ws.Range("I69").Value = 3562
$ws.Range("J69").Value = 3748.6667
$ws.Range("K69").Value = 10686
$ws.Range("L69").Value = 11246.0001
$ws.Range("M69").Value = -9812
$ws.Range("N69").Value = -12994.0001
$ws.Range("H72").Value = 3686.4443
$ws.Range("I72").Value = 3562
$ws.Range("J72").Value = 3748.6667
$ws.Range("K72").Value = 32058
$ws.Range("L72").Value = 33738.0003
$ws.Range("M72").Value = -27690
$ws.Range("N72").Value = -42474.0003
$ws.Range("H94").Value = 3699.8
$ws.Range("I94").Value = 3674.75
$ws.Range("J94").Value = 3800
$ws.Range("K94").Value = 3674.75
$ws.Range("L94").Value = 3800
$ws.Range("M94").Value = -3223.75
$ws.Range("N94").Value = -4702
$ws.Range("H111").Value = 3135.6667
$ws.Range("I111").Value = 5273
$ws.Range("J111").Value = 998.3333
$ws.Range("K111").Value = 15819
$ws.Range("L111").Value = 2994.9999
$ws.Range("M111").Value = -12752
$ws.Range("N111").Value = -9128.999899999999
$ws.Range("H137").Value = 1283.8
$ws.Range("I137").Value = 1299.8334
$ws.Range("J137").Value = 1259.75
$ws.Range("K137").Value = 3899.5002
$ws.Range("L137").Value = 3779.25
$ws.Range("M137").Value = -1349.5002
$ws.Range("N137").Value = -8879.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 315.33334
$ws.Range("I4").Value = 423
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 423
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = -307
$ws.Range("N4").Value = -332
$ws.Range("H32").Value = 11675.363
$ws.Range("I32").Value = 11675.363
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 11675.363
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -11388.363
$ws.Range("H61").Value = 1321.6786
$ws.Range("I61").Value = 956.2174
$ws.Range("J61").Value = 3002.8
$ws.Range("K61").Value = 956.2174
$ws.Range("L61").Value = 3002.8
$ws.Range("M61").Value = -744.2174
$ws.Range("N61").Value = -3426.8
$ws.Range("H74").Value = 1078.421
$ws.Range("I74").Value = 842.25
$ws.Range("J74").Value = 2338
$ws.Range("K74").Value = 842.25
$ws.Range("L74").Value = 2338
$ws.Range("M74").Value = 31.75
$ws.Range("N74").Value = -4086
$ws.Range("H77").Value = 1078.421
$ws.Range("I77").Value = 842.25
$ws.Range("J77").Value = 2338
$ws.Range("K77").Value = 4211.25
$ws.Range("L77").Value = 11690
$ws.Range("M77").Value = 156.75
$ws.Range("N77").Value = -20426
$ws.Range("H102").Value = 12821638
$ws.Range("I102").Value = 13890032
$ws.Range("J102").Value = 911
$ws.Range("K102").Value = 13890032
$ws.Range("L102").Value = 911
$ws.Range("M102").Value = -13888410
$ws.Range("N102").Value = -4155
$ws.Range("H107").Value = 17285.334
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 17285.334
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 17285.334
$ws.Range("N107").Value = -24965.334
$ws.Range("H136").Value = 1321.6786
$ws.Range("I136").Value = 956.2174
$ws.Range("J136").Value = 3002.8
$ws.Range("K136").Value = 2868.6522
$ws.Range("L136").Value = 9008.400000000001
$ws.Range("M136").Value = -318.6522
$ws.Range("N136").Value = -14108.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7350.8237
$ws.Range("I134").Value = 1464.3334
$ws.Range("J134").Value = 51499.5
$ws.Range("K134").Value = 4393.0002
$ws.Range("L134").Value = 154498.5
$ws.Range("M134").Value = -1858.0002
$ws.Range("N134").Value = -159568.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1738.25
$ws.Range("I31").Value = 865
$ws.Range("J31").Value = 2999.611
$ws.Range("K31").Value = 865
$ws.Range("L31").Value = 2999.611
$ws.Range("M31").Value = -570
$ws.Range("N31").Value = -3589.611
$ws.Range("H34").Value = 1738.25
$ws.Range("I34").Value = 865
$ws.Range("J34").Value = 2999.611
$ws.Range("K34").Value = 865
$ws.Range("L34").Value = 2999.611
$ws.Range("M34").Value = -663
$ws.Range("N34").Value = -3403.611
$ws.Range("H99").Value = 1723.1177
$ws.Range("I99").Value = 1807.375
$ws.Range("J99").Value = 1648.2222
$ws.Range("K99").Value = 1807.375
$ws.Range("L99").Value = 1648.2222
$ws.Range("M99").Value = -309.375
$ws.Range("N99").Value = -4644.2222
$ws.Range("H107").Value = 651.7273
$ws.Range("I107").Value = 470.91666
$ws.Range("J107").Value = 868.7
$ws.Range("K107").Value = 470.91666
$ws.Range("L107").Value = 868.7
$ws.Range("M107").Value = 1449.08334
$ws.Range("N107").Value = -4708.7
$ws.Range("H126").Value = 1723.1177
$ws.Range("I126").Value = 1807.375
$ws.Range("J126").Value = 1648.2222
$ws.Range("K126").Value = 5422.125
$ws.Range("L126").Value = 4944.6666
$ws.Range("M126").Value = -2952.125
$ws.Range("N126").Value = -9884.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1425.591
$ws.Range("I5").Value = 1510.65
$ws.Range("J5").Value = 575
$ws.Range("K5").Value = 4531.950000000001
$ws.Range("L5").Value = 1725
$ws.Range("M5").Value = -4419.950000000001
$ws.Range("N5").Value = -1949
$ws.Range("H80").Value = 4926.6665
$ws.Range("I80").Value = 2900
$ws.Range("J80").Value = 5071.4287
$ws.Range("K80").Value = 8700
$ws.Range("L80").Value = 15214.2861
$ws.Range("M80").Value = -7764
$ws.Range("N80").Value = -17086.2861
$ws.Range("H83").Value = 4926.6665
$ws.Range("I83").Value = 2900
$ws.Range("J83").Value = 5071.4287
$ws.Range("K83").Value = 26100
$ws.Range("L83").Value = 45642.85830000001
$ws.Range("M83").Value = -21420
$ws.Range("N83").Value = -55002.85830000001
$ws.Range("H92").Value = 800
$ws.Range("I92").Value = 400
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 1200
$ws.Range("L92").Value = 3000
$ws.Range("M92").Value = 48
$ws.Range("N92").Value = -5496
$ws.Range("H108").Value = 1424
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 1424
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 4272
$ws.Range("M108").ClearContents()
$ws.Range("N108").Value = -10032
$ws.Range("H113").Value = 748.1875
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 748.1875
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2244.5625
$ws.Range("N113").Value = -6584.5625
$ws.Range("H122").Value = 441.6
$ws.Range("I122").Value = 466
$ws.Range("J122").Value = 405
$ws.Range("K122").Value = 4194
$ws.Range("L122").Value = 3645
$ws.Range("M122").Value = -1744
$ws.Range("N122").Value = -8545
$ws.Range("H131").Value = 21740534
$ws.Range("I131").Value = 500000100
$ws.Range("J131").Value = 1461.7954
$ws.Range("K131").Value = 1500000300
$ws.Range("L131").Value = 4385.3862
$ws.Range("M131").Value = -1499995260
$ws.Range("N131").Value = -14465.3862
$ws.Range("H135").Value = 1425.591
$ws.Range("I135").Value = 1510.65
$ws.Range("J135").Value = 575
$ws.Range("K135").Value = 13595.85
$ws.Range("L135").Value = 5175
$ws.Range("M135").Value = -11060.85
$ws.Range("N135").Value = -10245
$ws.Range("H137").Value = 17579
$ws.Range("I137").Value = 3757.5
$ws.Range("J137").Value = 36007.668
$ws.Range("K137").Value = 11272.5
$ws.Range("L137").Value = 108023.004
$ws.Range("M137").Value = -6172.5
$ws.Range("N137").Value = -118223.004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2794.7778
$ws.Range("I113").Value = 1463.6666
$ws.Range("J113").Value = 4125.8887
$ws.Range("K113").Value = 1463.6666
$ws.Range("L113").Value = 4125.8887
$ws.Range("M113").Value = 706.3334
$ws.Range("N113").Value = -8465.8887

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2037.6666
$ws.Range("I7").Value = 1935.3334
$ws.Range("J7").Value = 2242.3333
$ws.Range("K7").Value = 1935.3334
$ws.Range("L7").Value = 2242.3333
$ws.Range("M7").Value = -1823.3334
$ws.Range("N7").Value = -2466.3333
$ws.Range("H40").Value = 6876.25
$ws.Range("I40").Value = 2250
$ws.Range("J40").Value = 11502.5
$ws.Range("K40").Value = 2250
$ws.Range("L40").Value = 11502.5
$ws.Range("M40").Value = -2114
$ws.Range("N40").Value = -11774.5
$ws.Range("H126").Value = 2037.6666
$ws.Range("I126").Value = 1935.3334
$ws.Range("J126").Value = 2242.3333
$ws.Range("K126").Value = 5806.0002
$ws.Range("L126").Value = 6726.999899999999
$ws.Range("M126").Value = -3336.0002
$ws.Range("N126").Value = -11666.9999
$ws.Range("H136").Value = 2090.2
$ws.Range("I136").Value = 1817.3334
$ws.Range("J136").Value = 2499.5
$ws.Range("K136").Value = 5452.0002
$ws.Range("L136").Value = 7498.5
$ws.Range("M136").Value = -2902.0002
$ws.Range("N136").Value = -12598.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 700005
$ws.Range("I14").Value = 700005
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 700005
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -699837
$ws.Range("H44").Value = 16280.5
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 16280.5
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 16280.5
$ws.Range("N44").Value = -17388.5
$ws.Range("H51").Value = 10038.25
$ws.Range("I51").Value = 9999
$ws.Range("J51").Value = 10051.333
$ws.Range("K51").Value = 9999
$ws.Range("L51").Value = 10051.333
$ws.Range("M51").Value = -9489
$ws.Range("N51").Value = -11071.333
